# Update custom sort order - restructure employee sheet with a new "Join" column
# inserted before "Position", fix a handful of data values, and append four new
# employee rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("employee")

# --- Step 1: restructure columns -----------------------------------------
# Old layout: A name | B gender | C DoB | D Position | E Salary | F Join
# New layout: A name | B gender | C DoB | D Join | E Position | F Salary
# Moving the "Join" column (F) so it sits right after DoB (before Position)
# is equivalent to cutting column F and inserting it at column D.
$ws.Columns("F").Cut()
$ws.Columns("D").Insert()

# --- Step 2: column B (gender) got a little wider -------------------------
# (target stored width is 10.21875 characters; the engine quantizes
# ColumnWidth to whole-pixel steps, so 9.25 is the closest settable value
# that lands on the nearest achievable stored width)
$ws.Columns("B").ColumnWidth = 9.25

# --- Step 3: fix a handful of data values that changed alongside the
# column restructuring (values below use the NEW column layout) ----------

# Row 2 : Join date corrected
$ws.Range("D2").Value = 35723

# Row 5 : DoB and Join dates corrected
$ws.Range("C5").Value = 32350
$ws.Range("D5").Value = 42702

# Row 10 : DoB corrected
$ws.Range("C10").Value = 35926

# Row 14 : DoB and Join dates corrected
$ws.Range("C14").Value = 33598
$ws.Range("D14").Value = 41314

# Row 15 : Salary corrected
$ws.Range("F15").Value = 22000

# Row 20 : Salary corrected
$ws.Range("F20").Value = 77000

# Row 21 : Join date, Position and Salary corrected
$ws.Range("D21").Value = 36349
$ws.Range("E21").Value = "ผู้จัดการ"
$ws.Range("F21").Value = 90000

# --- Step 4: append four new employee rows --------------------------------
$ws.Range("A22").Value = "ฟ้า"
$ws.Range("B22").Value = "F"
$ws.Range("C22").Value = 38222
$ws.Range("D22").Value = 45047
$ws.Range("E22").Value = "พนักงานระดับต้น"
$ws.Range("F22").Value = 18500

$ws.Range("A23").Value = "อิงฟ้า"
$ws.Range("B23").Value = "F"
$ws.Range("C23").Value = 37869
$ws.Range("D23").Value = 45078
$ws.Range("E23").Value = "พนักงานระดับต้น"
$ws.Range("F23").Value = 19000

$ws.Range("A24").Value = "สายฟ้า"
$ws.Range("B24").Value = "M"
$ws.Range("C24").Value = 37655
$ws.Range("D24").Value = 45108
$ws.Range("E24").Value = "พนักงานระดับต้น"
$ws.Range("F24").Value = 19500

$ws.Range("A25").Value = "สายรุ้ง"
$ws.Range("B25").Value = "F"
$ws.Range("C25").Value = 37389
$ws.Range("D25").Value = 45139
$ws.Range("E25").Value = "พนักงานระดับต้น"
$ws.Range("F25").Value = 20000

# Apply the date/currency number formats used by the rest of the table to
# the newly appended rows (matching columns C/D date style, F currency style)
# by copying formatting from the row above, so we reuse the existing style
# entries rather than minting new ones.
$ws.Range("C21:D21").Copy()
$ws.Range("C22:D25").PasteSpecial(-4122)
$ws.Range("F21").Copy()
$ws.Range("F22:F25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Step 5: update the active selection on the sheet ----------------------
$ws.Range("C10").Select()

"done"
